$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.420398
$ws.Range("H2").Value = 1.261194
$ws.Range("I2").Value = 0.001794717364332138
$ws.Range("J2").Value = 0.001794717364332138
$ws.Range("M2").Value = 1.365066333333333
$ws.Range("N2").Value = 4.095199
$ws.Range("O2").Value = 0.1025111867562684
$ws.Range("P2").Value = 0.1025111867562683
$ws.Range("Q2").Value = 0.5738711564006667
$ws.Range("R2").Value = 5.164840407606
$ws.Range("S2").Value = 0.0001839786069097695
$ws.Range("T2").Value = 0.0001839786069097695

$ws.Range("G3").Value = 0.420398
$ws.Range("H3").Value = 1.261194
$ws.Range("I3").Value = 0.001794717364332138
$ws.Range("J3").Value = 0.001794717364332138
$ws.Range("M3").Value = 5.511188000000001
$ws.Range("N3").Value = 16.533564
$ws.Range("O3").Value = 0.4138688417707456
$ws.Range("P3").Value = 0.4138688417707455
$ws.Range("Q3").Value = 2.316892412824
$ws.Range("R3").Value = 20.852031715416
$ws.Range("S3").Value = 0.000742777596881987
$ws.Range("T3").Value = 0.0007427775968819869

$ws.Range("G4").Value = 0.420398
$ws.Range("H4").Value = 1.261194
$ws.Range("I4").Value = 0.001794717364332138
$ws.Range("J4").Value = 0.001794717364332138
$ws.Range("M4").Value = 4.174538333333333
$ws.Range("N4").Value = 12.523615
$ws.Range("O4").Value = 0.3134916364573745
$ws.Range("P4").Value = 0.3134916364573745
$ws.Range("Q4").Value = 1.754967566256667
$ws.Range("R4").Value = 15.79470809631
$ws.Range("S4").Value = 0.000562628883522948
$ws.Range("T4").Value = 0.0005626288835229479

$ws.Range("G5").Value = 0.420398
$ws.Range("H5").Value = 1.261194
$ws.Range("I5").Value = 0.001794717364332138
$ws.Range("J5").Value = 0.001794717364332138
$ws.Range("M5").Value = 2.265474333333334
$ws.Range("N5").Value = 6.796423000000001
$ws.Range("O5").Value = 0.1701283350156116
$ws.Range("P5").Value = 0.1701283350156116
$ws.Range("Q5").Value = 0.9524008787846668
$ws.Range("R5").Value = 8.571607909062003
$ws.Range("S5").Value = 0.0003053322770174334
$ws.Range("T5").Value = 0.0003053322770174334

$ws.Range("G6").Value = 0.04680500000000001
$ws.Range("H6").Value = 0.140415
$ws.Range("I6").Value = 0.0001998148093891163
$ws.Range("J6").Value = 0.0001998148093891163
$ws.Range("M6").Value = 1.365066333333333
$ws.Range("N6").Value = 4.095199
$ws.Range("O6").Value = 0.1025111867562684
$ws.Range("P6").Value = 0.1025111867562683
$ws.Range("Q6").Value = 0.06389192973166667
$ws.Range("R6").Value = 0.5750273675850001
$ws.Range("S6").Value = 0.00002048325324195586
$ws.Range("T6").Value = 0.00002048325324195586

$ws.Range("G7").Value = 0.04680500000000001
$ws.Range("H7").Value = 0.140415
$ws.Range("I7").Value = 0.0001998148093891163
$ws.Range("J7").Value = 0.0001998148093891163
$ws.Range("M7").Value = 5.511188000000001
$ws.Range("N7").Value = 16.533564
$ws.Range("O7").Value = 0.4138688417707456
$ws.Range("P7").Value = 0.4138688417707455
$ws.Range("Q7").Value = 0.25795115434
$ws.Range("R7").Value = 2.321560389060001
$ws.Range("S7").Value = 0.00008269712373051585
$ws.Range("T7").Value = 0.00008269712373051583

$ws.Range("G8").Value = 0.04680500000000001
$ws.Range("H8").Value = 0.140415
$ws.Range("I8").Value = 0.0001998148093891163
$ws.Range("J8").Value = 0.0001998148093891163
$ws.Range("M8").Value = 4.174538333333333
$ws.Range("N8").Value = 12.523615
$ws.Range("O8").Value = 0.3134916364573745
$ws.Range("P8").Value = 0.3134916364573745
$ws.Range("Q8").Value = 0.1953892666916667
$ws.Range("R8").Value = 1.758503400225
$ws.Range("S8").Value = 0.00006264027158381243
$ws.Range("T8").Value = 0.00006264027158381242

$ws.Range("G9").Value = 0.04680500000000001
$ws.Range("H9").Value = 0.140415
$ws.Range("I9").Value = 0.0001998148093891163
$ws.Range("J9").Value = 0.0001998148093891163
$ws.Range("M9").Value = 2.265474333333334
$ws.Range("N9").Value = 6.796423000000001
$ws.Range("O9").Value = 0.1701283350156116
$ws.Range("P9").Value = 0.1701283350156116
$ws.Range("Q9").Value = 0.1060355261716667
$ws.Range("R9").Value = 0.9543197355450002
$ws.Range("S9").Value = 0.00003399416083283215
$ws.Range("T9").Value = 0.00003399416083283215

$ws.Range("G10").Value = 233.774694
$ws.Range("H10").Value = 701.324082
$ws.Range("I10").Value = 0.9980054678262787
$ws.Range("J10").Value = 0.9980054678262787
$ws.Range("M10").Value = 1.365066333333333
$ws.Range("N10").Value = 4.095199
$ws.Range("O10").Value = 0.1025111867562684
$ws.Range("P10").Value = 0.1025111867562683
$ws.Range("Q10").Value = 319.1179643647019
$ws.Range("R10").Value = 2872.061679282318
$ws.Range("S10").Value = 0.1023067248961166
$ws.Range("T10").Value = 0.1023067248961166

$ws.Range("G11").Value = 233.774694
$ws.Range("H11").Value = 701.324082
$ws.Range("I11").Value = 0.9980054678262787
$ws.Range("J11").Value = 0.9980054678262787
$ws.Range("M11").Value = 5.511188000000001
$ws.Range("N11").Value = 16.533564
$ws.Range("O11").Value = 0.4138688417707456
$ws.Range("P11").Value = 0.4138688417707455
$ws.Range("Q11").Value = 1288.376288276472
$ws.Range("R11").Value = 11595.38659448825
$ws.Range("S11").Value = 0.4130433670501331
$ws.Range("T11").Value = 0.413043367050133

$ws.Range("G12").Value = 233.774694
$ws.Range("H12").Value = 701.324082
$ws.Range("I12").Value = 0.9980054678262787
$ws.Range("J12").Value = 0.9980054678262787
$ws.Range("M12").Value = 4.174538333333333
$ws.Range("N12").Value = 12.523615
$ws.Range("O12").Value = 0.3134916364573745
$ws.Range("P12").Value = 0.3134916364573745
$ws.Range("Q12").Value = 975.90142146627
$ws.Range("R12").Value = 8783.11279319643
$ws.Range("S12").Value = 0.3128663673022678
$ws.Range("T12").Value = 0.3128663673022677

$ws.Range("G13").Value = 233.774694
$ws.Range("H13").Value = 701.324082
$ws.Range("I13").Value = 0.9980054678262787
$ws.Range("J13").Value = 0.9980054678262787
$ws.Range("M13").Value = 2.265474333333334
$ws.Range("N13").Value = 6.796423000000001
$ws.Range("O13").Value = 0.1701283350156116
$ws.Range("P13").Value = 0.1701283350156116
$ws.Range("Q13").Value = 529.610569039854
$ws.Range("R13").Value = 4766.495121358686
$ws.Range("S13").Value = 0.1697890085777614
$ws.Range("T13").Value = 0.1697890085777614

